$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert 5 new rows before row 520, pushing the existing
# rows 520-580 down to 525-585 (new used range A1:T585), then populate the
# 5 freshly-inserted rows with this week's price quotes.
$ws.Rows("520:524").Insert()

$row520 = New-Object "object[,]" 1,20
$row520[0,0] = 5
$row520[0,1] = 'Macroferia Regional de Talca'
$row520[0,2] = 'Maule'
$row520[0,3] = 44474
$row520[0,4] = 7
$row520[0,5] = 'Fruta'
$row520[0,6] = 100106
$row520[0,7] = 'Oleaginosos'
$row520[0,8] = 100106002
$row520[0,9] = 'Palta'
$row520[0,10] = 'Hass'
$row520[0,11] = '1a nueva(o)'
$row520[0,12] = 100
$row520[0,13] = 2600
$row520[0,14] = 2600
$row520[0,15] = 2600
$row520[0,16] = '$/kilo (en caja de 17 kilos)'
$row520[0,17] = 'Cabildo'
$row520[0,18] = 2600
$row520[0,19] = 1
$ws.Range("A520:T520").Value = $row520

$row521 = New-Object "object[,]" 1,20
$row521[0,0] = 5
$row521[0,1] = 'Macroferia Regional de Talca'
$row521[0,2] = 'Maule'
$row521[0,3] = 44474
$row521[0,4] = 7
$row521[0,5] = 'Fruta'
$row521[0,6] = 100106
$row521[0,7] = 'Oleaginosos'
$row521[0,8] = 100106002
$row521[0,9] = 'Palta'
$row521[0,10] = 'Hass'
$row521[0,11] = '2a nueva(o)'
$row521[0,12] = 130
$row521[0,13] = 2300
$row521[0,14] = 2300
$row521[0,15] = 2300
$row521[0,16] = '$/kilo (en caja de 17 kilos)'
$row521[0,17] = 'Cabildo'
$row521[0,18] = 2300
$row521[0,19] = 1
$ws.Range("A521:T521").Value = $row521

$row522 = New-Object "object[,]" 1,20
$row522[0,0] = 5
$row522[0,1] = 'Macroferia Regional de Talca'
$row522[0,2] = 'Maule'
$row522[0,3] = 44474
$row522[0,4] = 7
$row522[0,5] = 'Fruta'
$row522[0,6] = 100106
$row522[0,7] = 'Oleaginosos'
$row522[0,8] = 100106002
$row522[0,9] = 'Palta'
$row522[0,10] = 'Hass'
$row522[0,11] = 'Especial nueva (o)'
$row522[0,12] = 250
$row522[0,13] = 3000
$row522[0,14] = 3000
$row522[0,15] = 3000
$row522[0,16] = '$/kilo (en caja de 17 kilos)'
$row522[0,17] = 'Cabildo'
$row522[0,18] = 3000
$row522[0,19] = 1
$ws.Range("A522:T522").Value = $row522

$row523 = New-Object "object[,]" 1,20
$row523[0,0] = 5
$row523[0,1] = 'Macroferia Regional de Talca'
$row523[0,2] = 'Maule'
$row523[0,3] = 44474
$row523[0,4] = 7
$row523[0,5] = 'Fruta'
$row523[0,6] = 100106
$row523[0,7] = 'Oleaginosos'
$row523[0,8] = 100106002
$row523[0,9] = 'Palta'
$row523[0,10] = 'Negra de La Cruz'
$row523[0,11] = 'Especial'
$row523[0,12] = 200
$row523[0,13] = 2500
$row523[0,14] = 2500
$row523[0,15] = 2500
$row523[0,16] = '$/kilo (en caja de 8 kilos )'
$row523[0,17] = 'Provincia de Melipilla'
$row523[0,18] = 2500
$row523[0,19] = 1
$ws.Range("A523:T523").Value = $row523

$row524 = New-Object "object[,]" 1,20
$row524[0,0] = 5
$row524[0,1] = 'Macroferia Regional de Talca'
$row524[0,2] = 'Maule'
$row524[0,3] = 44474
$row524[0,4] = 7
$row524[0,5] = 'Fruta'
$row524[0,6] = 100106
$row524[0,7] = 'Oleaginosos'
$row524[0,8] = 100106002
$row524[0,9] = 'Palta'
$row524[0,10] = 'Negra de La Cruz'
$row524[0,11] = 'Primera'
$row524[0,12] = 250
$row524[0,13] = 2000
$row524[0,14] = 2000
$row524[0,15] = 2000
$row524[0,16] = '$/kilo (en caja de 8 kilos )'
$row524[0,17] = 'Provincia de Melipilla'
$row524[0,18] = 2000
$row524[0,19] = 1
$ws.Range("A524:T524").Value = $row524

